$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.389.23"
$ws.Range("E2").Value = "  +3.91%  "

$ws.Range("D3").Value = "3.084.23"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").Value = "'217.96"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").Value = "'616.95"
$ws.Range("E6").Value = "  -2.75%  "

$ws.Range("E7").Value = "  -3.36%  "

$ws.Range("D8").Value = "'0.901"
$ws.Range("E8").Value = "  +8.00%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "3.081.64"
$ws.Range("E10").Value = "  -0.93%  "

$ws.Range("D11").Value = "'0.670"
$ws.Range("E11").Value = "  +16.41%  "

$ws.Range("E12").Value = "  +6.00%  "

$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").Value = "91.171.25"
$ws.Range("E14").Value = "  +3.67%  "

$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "'32.98"
$ws.Range("E16").Value = "  +3.11%  "

$ws.Range("D17").Value = "3.658.25"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "3.077.70"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").Value = "'3.46"
$ws.Range("E19").Value = "  +2.70%  "

$ws.Range("D20").Value = "'0.0000221"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "'13.79"
$ws.Range("E21").Value = "  +4.36%  "

$ws.Range("D22").Value = "'434.92"
$ws.Range("E22").Value = "  +2.62%  "

$ws.Range("D23").Value = "'8.46"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("D24").Value = "'5.14"
$ws.Range("E24").Value = "  +5.14%  "

$ws.Range("D25").Value = "'5.64"
$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").Value = "'84.00"
$ws.Range("E26").Value = "  +1.20%  "

$ws.Range("D27").Value = "'11.78"
$ws.Range("E27").Value = "  +2.53%  "

$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D30").Value = "'0.167"
$ws.Range("E30").Value = "  +7.00%  "

$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("D32").Value = "'8.79"
$ws.Range("E32").Value = "  +7.84%  "

$ws.Range("D33").Value = "'3.88"
$ws.Range("E33").Value = "  -2.11%  "

$ws.Range("D34").Value = "'518.50"
$ws.Range("E34").Value = "  +3.80%  "

$ws.Range("D35").Value = "'7.05"
$ws.Range("E35").Value = "  +4.02%  "

$ws.Range("E36").Value = "  -7.07%  "

$ws.Range("D37").Value = "'1.28"
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("E39").Value = "  +2.98%  "

$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E43").Value = "  +1.90%  "

$ws.Range("D44").Value = "'0.368"
$ws.Range("E44").Value = "  +0.74%  "

$ws.Range("E45").Value = "  +1.83%  "

$ws.Range("E46").Value = "  +10.12%  "

$ws.Range("D47").Value = "'43.89"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").Value = "'141.53"
$ws.Range("E48").Value = "  -3.08%  "

$ws.Range("D49").Value = "'0.000264"
$ws.Range("E49").Value = "  +12.99%  "

$ws.Range("E50").Value = "  +6.63%  "

$ws.Range("D51").Value = "'164.50"
$ws.Range("E51").Value = "  +1.46%  "
